$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 449, shifting existing rows 449:473 down to 450:474
$ws.Rows.Item(449).Insert()

# Populate the newly inserted row 449 with the new weekly price entry
$ws.Cells.Item(449, 1).Value = 10
$ws.Cells.Item(449, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(449, 3).Value = "La Araucanía"
$ws.Cells.Item(449, 4).Value = 45267
$ws.Cells.Item(449, 5).Value = 9
$ws.Cells.Item(449, 6).Value = "Fruta"
$ws.Cells.Item(449, 7).Value = 100103
$ws.Cells.Item(449, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(449, 9).Value = 100103004
$ws.Cells.Item(449, 10).Value = "Durazno"
$ws.Cells.Item(449, 11).Value = "Early Majestic"
$ws.Cells.Item(449, 12).Value = "Primera"
$ws.Cells.Item(449, 13).Value = 250
$ws.Cells.Item(449, 14).Value = 20000
$ws.Cells.Item(449, 15).Value = 20000
$ws.Cells.Item(449, 16).Value = 20000
$ws.Cells.Item(449, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(449, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(449, 19).Value = 1111
$ws.Cells.Item(449, 20).Value = 18
